$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing content (and the
# "namespace" header row + hyperlink) down by one row.
$ws.Rows("1:1").Insert()

# New title row
$ws.Range("A1").Value = "title"
$ws.Range("A1").Font.Bold = $true

$ws.Range("B1").Value = "Campos de la Investigación y el desarrollo (OCDE)"

# Row insertion does not relocate the hyperlink annotation that used to
# live on the old B1 ("namespace" URL); move it onto the new B2.
$ws.Range("B1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "http://purl.org/concytec-pe/ford_ocde")

# Update selection to match target workbook state
$ws.Range("B9").Select()
